$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("OrangeHRM_AddUser")

# --- Row 2 edits: shorten the Charlie Carter name ---
$ws.Range("F2").Value() = "Charlie Car"

# --- New data for rows 3-11, columns F (Sarika_n) and G (Charlie_n) ---
$sarika = @("Sarika_1","Sarika_2","Sarika_3","Sarika_4","Sarika_5","Sarika_6","Sarika_7","Sarika_8","Sarika_9")
$charlie = @("Charlie_1","Charlie_2","Charlie_3","Charlie_4","Charlie_5","Charlie_6","Charlie_7","Charlie_8","Charlie_9")
$admin = @("Admin@123","Admin@124","Admin@125","Admin@126","Admin@127","Admin@128","Admin@129","Admin@130","Admin@131")

for ($i = 0; $i -lt 9; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 6).Value() = $sarika[$i]
    $ws.Cells.Item($row, 7).Value() = $charlie[$i]
}

# --- Column H (8) hyperlinked Admin@ values ---
# H3 gets its own single-cell hyperlink (no display override)
$h3 = $ws.Range("H3")
$ws.Hyperlinks.Add($h3, "mailto:Admin@123") | Out-Null
$h3.Value() = "Admin@123"
$h3.Style = $ws.Range("H2").Style

# H4:H11 share one hyperlink with display text "Admin@123", but each cell keeps its own value
$hRange = $ws.Range("H4:H11")
$ws.Hyperlinks.Add($hRange, "mailto:Admin@123", "", "", "Admin@123") | Out-Null
for ($i = 1; $i -lt 9; $i++) {
    $row = 3 + $i
    $cell = $ws.Cells.Item($row, 8)
    $cell.Value() = $admin[$i]
    $cell.Style = $ws.Range("H2").Style
}

# --- Column I (9) mirrors column H ---
$i3 = $ws.Range("I3")
$ws.Hyperlinks.Add($i3, "mailto:Admin@123") | Out-Null
$i3.Value() = "Admin@123"
$i3.Style = $ws.Range("I2").Style

$iRange = $ws.Range("I4:I11")
$ws.Hyperlinks.Add($iRange, "mailto:Admin@123", "", "", "Admin@123") | Out-Null
for ($i = 1; $i -lt 9; $i++) {
    $row = 3 + $i
    $cell = $ws.Cells.Item($row, 9)
    $cell.Value() = $admin[$i]
    $cell.Style = $ws.Range("I2").Style
}

# --- Column H width set manually (renders as width=18 in the saved file) ---
$ws.Columns.Item(8).ColumnWidth = 17.1

# --- Final selection left on I14 ---
$ws.Activate() | Out-Null
$ws.Range("I14").Select() | Out-Null
